{"js": "// Replace the Slovak \"V roku 2018 ... Perseus ...\" observation-dates\n// sentence with the new Bootes observation dates, as a single plain run\n// (no direct character formatting), in every paragraph that contains it.\n\nconst oldMarker = \"V roku 2018\";\nconst newText = \"V roku Bootes: 14.-23. m\u00e1ja, 13.-22. j\u00fana, 12.-21. j\u00fala\";\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst targets = [];\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text.indexOf(oldMarker) !== -1) {\n    targets.push(paragraphs.items[i]);\n  }\n}\n\nfor (const paragraph of targets) {\n  const range = paragraph.getRange();\n  // Clear() removes both the paragraph's text AND its direct (run-level)\n  // character formatting, leaving the paragraph mark/pPr untouched.\n  range.clear();\n  // Insert the new sentence into the now-empty, unformatted paragraph so\n  // it lands as a single run with no rPr, matching the target edit.\n  range.insertText(newText, \"Start\");\n}\n\nawait context.sync();\n", "ps1": "# Replace the Slovak \"V roku 2018 ... Perseus ...\" observation-dates\n# sentence with the new Bootes observation dates, as a single plain run\n# (no direct character formatting), in every paragraph that contains it.\n\n$d = $word.ActiveDocument\n$oldMarker = \"V roku 2018\"\n$newText = \"V roku Bootes: 14.-23. m\u00e1ja, 13.-22. j\u00fana, 12.-21. j\u00fala\"\n\n# Snapshot the paragraph ranges first: mutating paragraph text while\n# iterating $d.Paragraphs directly can disturb the live enumerator.\n$targets = @()\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text.Contains($oldMarker)) {\n        $targets += $p.Range\n    }\n}\n\nforeach ($r in $targets) {\n    # $r currently spans the whole paragraph, including its trailing\n    # paragraph mark; back off one character so we only remove the\n    # paragraph's text content (runs) and keep the paragraph itself.\n    $r.MoveEnd(1, -1) | Out-Null   # wdCharacter = 1\n    $r.Delete()\n    # Insert fresh text into the now-empty paragraph; with no preceding\n    # character formatting to inherit, this lands as a single run with\n    # no rPr (direct formatting), matching the target edit.\n    $r.InsertAfter($newText)\n}\n"}
